# Update stats for 2026-01
# Append a new data row (row 26) to Sheet1 mirroring the existing
# month/schools/authorities/users/... columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 26
$prevRow = $newRow - 1

$ws.Cells.Item($newRow, 1).Value = 46023
$ws.Cells.Item($newRow, 2).Value = 6479
$ws.Cells.Item($newRow, 3).Value = 1010
$ws.Cells.Item($newRow, 4).Value = 6043294
$ws.Cells.Item($newRow, 5).Value = 932.7510418274425
$ws.Cells.Item($newRow, 6).Value = 9.498056447524084
$ws.Cells.Item($newRow, 7).Value = 7.218683651804669
$ws.Cells.Item($newRow, 8).Value = 25.8518011467513

# Match the date/number formatting used by the month column above it.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat
